$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column D ("NOTA") header + widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 94.5703125
$ws.Columns.Item(4).ColumnWidth = 34.7109375

$ws.Range("D1").Value = "NOTA"
$ws.Range("D1").Font.Name = "Arial"
$ws.Range("D1").Font.Size = 11
$ws.Range("D1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Row 2 (existing "consultarTodosVotantes" row) gets a real hyperlink + the
# smaller Arial 10 font used by the new rows; C2 wraps, D2 stays blank but
# picks up the same Arial 10 font as the rest of the row.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A2"), "http://192.190.43.117/VotoBack/v1/VotantesService/consultarTodosVotantes") | Out-Null

$ws.Range("B2:D2").Font.Name = "Arial"
$ws.Range("B2:D2").Font.Size = 10
$ws.Range("C2").WrapText = $true

# ---------------------------------------------------------------------------
# Row 3: "validarVotante" service
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "http://192.190.43.117/VotoBack/v1/VotantesService/validarVotante/{cedula}"
$ws.Hyperlinks.Add($ws.Range("A3"), "http://192.190.43.117/VotoBack/v1/VotantesService/validarVotante/{cedula}") | Out-Null

$ws.Range("B3").Value = "GET"

$ws.Range("C3").Value = "VALIDA SI UN VOTANTE:`n-EXISTE`n-YA ESTA VOTANDO`n-YA VOTÓ"
$ws.Range("D3").Value = "CUANDO SE INVOQUE ESTE SERVICIO SE ENTENDERA QUE EL VOTANTE PROCEDERA A VOTAR"

$ws.Range("B3:D3").Font.Name = "Arial"
$ws.Range("B3:D3").Font.Size = 10
$ws.Range("C3:D3").WrapText = $true

$ws.Rows.Item(3).RowHeight = 51.75

# ---------------------------------------------------------------------------
# Row 4: "guardarVoto" service
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "http://192.190.43.117/VotoBack/v1/VotantesService/guardarVoto/{cedula}/{idCandidato}/{idPuesto}"
$ws.Hyperlinks.Add($ws.Range("A4"), "http://192.190.43.117/VotoBack/v1/VotantesService/guardarVoto/{cedula}/{idCandidato}/{idPuesto}") | Out-Null

$ws.Range("B4").Value = "POST"
$ws.Range("C4").Value = "GUARDA UN VOTO"

# ---------------------------------------------------------------------------
# Final selection, matching the authored file
# ---------------------------------------------------------------------------
$ws.Range("C4").Select()
